$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition) updates to column F ("想去人数" / "want-to-go count")
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 3362
$wsExpo.Range("F6").Value = 2355
$wsExpo.Range("F8").Value = 99
$wsExpo.Range("F13").Value = 171
$wsExpo.Range("F14").Value = 561

# Sheet "全部类型" (all types) updates to column F ("想去人数" / "want-to-go count")
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 3362
$wsAll.Range("F7").Value = 2355
$wsAll.Range("F9").Value = 99
$wsAll.Range("F14").Value = 171
$wsAll.Range("F15").Value = 561
